$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Fcaliper_offset (mm)"
$ws.Range("B18").Value = 0

$ws.Range("A19").Value = "Fcaliper_radius (mm)"
$ws.Range("B19").Value = 90

$ws.Range("A20").Value = "Fcaliper_distance (mm)"
$ws.Range("B20").Value = 69.85

$ws.Range("B21").Select()
